# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Agosto de 2020 a las 18:12"

# --- Rows whose country identity + stats swap places (shared-string reorder) ---
# Jordania / Malta swap (rows 144 / 145)
$ws.Range("A144").Value = "Jordania"
$ws.Range("B144").Value = 1438
$ws.Range("C144").Value = 40
$ws.Range("D144").Value = 1241
$ws.Range("E144").Value = 186
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 11

$ws.Range("A145").Value = "Malta"
$ws.Range("B145").Value = 1423
$ws.Range("C145").Value = 48
$ws.Range("D145").Value = 766
$ws.Range("E145").Value = 648
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 9

# Islas Feroe / Birmania swap (rows 173 / 174)
$ws.Range("A173").Value = "Islas Feroe"
$ws.Range("B173").Value = 377
$ws.Range("C173").Value = 4
$ws.Range("D173").Value = 235
$ws.Range("E173").Value = 142
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 0

$ws.Range("A174").Value = "Birmania"
$ws.Range("B174").Value = 376
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 331
$ws.Range("E174").Value = 39
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 6

# Islas Malvinas / Montserrat swap (rows 213 / 214)
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 13
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1

# --- Standalone statistic refreshes (country identity unchanged) ---
# Row 4: Estados Unidos
$ws.Range("B4").Value = 5622237
$ws.Range("C4").Value = 10210
$ws.Range("D4").Value = 2974803
$ws.Range("E4").Value = 2473328
$ws.Range("G4").Value = 390
$ws.Range("H4").Value = 174106

# Row 15: Reino Unido
$ws.Range("B15").Value = 320286
$ws.Range("C15").Value = 1089
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = 41381

# Row 17: Argentina
$ws.Range("D17").Value = 223531
$ws.Range("E17").Value = 69718
$ws.Range("G17").Value = 63
$ws.Range("H17").Value = 5877

# Row 20: Italia
$ws.Range("B20").Value = 254636
$ws.Range("C20").Value = 401
$ws.Range("D20").Value = 204142
$ws.Range("E20").Value = 15089
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 35405

# Row 48: Singapur
$ws.Range("D48").Value = 52533
$ws.Range("E48").Value = 3378

# Row 53: Barein
$ws.Range("E53").Value = 3481
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 175

# Row 63: Azerbaiyan
$ws.Range("B63").Value = 34474
$ws.Range("C63").Value = 131
$ws.Range("D63").Value = 32201
$ws.Range("E63").Value = 1764
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 509

# Row 135: Islandia
$ws.Range("B135").Value = 2027
$ws.Range("C135").Value = 13
$ws.Range("D135").Value = 1895
$ws.Range("E135").Value = 122
